{"js": "/*\n * The document has a stray \"_GoBack\" bookmark \u2014 an empty\n * bookmarkStart/bookmarkEnd pair that Word silently drops on the\n * last-edit-position whenever it next saves the file. Clean it up here\n * to mirror that behavior (this also matches the namespace/compat\n * bump recorded for this save in the source diff).\n */\n\n// Word.Document has no public getBookmarks(); Body/Range expose it, and\n// \"_GoBack\" is a hidden bookmark, so ask for hidden ones explicitly.\nconst bookmarks = context.document.body.getBookmarks(true /* includeHidden */, false /* includeAdjacent */);\nawait context.sync();\n\nif (bookmarks.value && bookmarks.value.indexOf(\"_GoBack\") !== -1) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The document has a stray \"_GoBack\" bookmark - an empty\n# bookmarkStart/bookmarkEnd pair that Word silently drops on the\n# last-edit-position whenever it next saves the file. Clean it up here\n# to mirror that behavior (this also matches the namespace/compat\n# bump recorded for this save in the source diff).\n\n$d = $word.ActiveDocument\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
